$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J, matching style of existing header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J21
$data = @(
    @(1, 5),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(1, 4),
    @(1, 4),
    @(1, 3),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(5, 6),
    @(7, 7),
    @(5, 6),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(4, 5),
    @(5, 6),
    @(7, 7),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
